$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 values to the rounded (2-decimal) "custom accuracy" figures.
$ws.Range("B5").Value = 21.14
$ws.Range("C5").Value = 15.54
$ws.Range("D5").Value = 1.26
$ws.Range("E5").Value = 45.87
$ws.Range("F5").Value = 37.64
$ws.Range("G5").Value = 16.64
$ws.Range("H5").Value = 65.24
$ws.Range("I5").Value = 25.6
$ws.Range("J5").Value = 11.34
$ws.Range("K5").Value = 16.83
$ws.Range("L5").Value = 18.43
$ws.Range("M5").Value = 19.39
$ws.Range("N5").Value = 5.31
$ws.Range("O5").Value = 16.54
$ws.Range("P5").Value = 23.53
$ws.Range("Q5").Value = 13.96
$ws.Range("R5").Value = 0.86
$ws.Range("S5").Value = 0.86
$ws.Range("T5").Value = 244.48
$ws.Range("U5").Value = 46.29
$ws.Range("V5").Value = 15.27
$ws.Range("W5").Value = 31.1
$ws.Range("X5").Value = 16.35
$ws.Range("Y5").Value = 2.48
$ws.Range("Z5").Value = 31.46
$ws.Range("AA5").Value = 13.49
$ws.Range("AB5").Value = 11.98
$ws.Range("AC5").Value = 14.08
$ws.Range("AD5").Value = 19.27
$ws.Range("AE5").Value = 0.56
$ws.Range("AF5").Value = 59.22
$ws.Range("AG5").Value = 8.6
$ws.Range("AH5").Value = 19.09

# Remove the now-superfluous last data row (row 6); this also shrinks the
# sheet dimension from A1:AH6 down to A1:AH5.
$ws.Rows.Item(6).Delete()
